# -- added ROC curves to RF report --
# Bump the cached "as of" date shown in the poster footer (master + every
# layout carries its own copy of the date placeholder), and refresh the
# Random forest "Good" row Sens. values in the classifier-performance table
# now that ROC curves were added to the report.

function Update-DatePlaceholder($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.HasTextFrame) {
                $shp.TextFrame.TextRange.Text = $newText
            }
        }
    }
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster

Update-DatePlaceholder $master.Shapes "11/04/2014"

$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    Update-DatePlaceholder $layouts.Item($L).Shapes "11/04/2014"
}

# Table 10 ("Random forest" / "Support vector machine" classifier report),
# row "Good": Sens. columns for the Random forest model move from 0.85 to
# 0.88 (Linear model) and 0.82 (Correlation).
$s = $p.Slides.Item(1)
$tbl = $s.Shapes.Item(1).Table
$tbl.Cell(6, 2).Shape.TextFrame.TextRange.Text = "0.88"
$tbl.Cell(6, 4).Shape.TextFrame.TextRange.Text = "0.82"
